# revision 8: major changes to text and figures, completed supplementary
# material, almost ready for submission.
#
# 1) Bump the "update automatically" date fields from 7/10/2013 to
#    7/17/2013 on the slide master and every slide layout.
# 2) Merge the two adjacent 8-space runs at the top of the "TextBox 11"
#    shape into a single 16-space run.
# 3) Add "c)" / "d)" sub-figure labels under the existing "a)" / "b)"
#    labels, modelled on the existing "b)" textbox.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholders: slide master + every custom (slide) layout.
#    (NotesMaster is intentionally skipped -- writing through
#    $p.NotesMaster.Shapes in this host mis-routes to the SlideMaster,
#    so touching it would corrupt unrelated master content.)
# ---------------------------------------------------------------------
$masterShapes = $p.SlideMaster.Shapes
for ($i = 1; $i -le $masterShapes.Count; $i++) {
    $shp = $masterShapes.Item($i)
    if ($shp.Name -like "*Date*") {
        if ($shp.TextFrame.TextRange.Text -eq "7/10/2013") {
            $shp.TextFrame.TextRange.Text = "7/17/2013"
        }
    }
}

$layouts = $p.Designs.Item(1).SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    $layoutShapes = $layout.Shapes
    for ($j = 1; $j -le $layoutShapes.Count; $j++) {
        $lshp = $layoutShapes.Item($j)
        if ($lshp.Name -like "*Date*") {
            if ($lshp.TextFrame.TextRange.Text -eq "7/10/2013") {
                $lshp.TextFrame.TextRange.Text = "7/17/2013"
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2) Collapse the two 8-space runs into one 16-space run on slide 1.
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(1)
$slideShapes = $slide.Shapes
for ($k = 1; $k -le $slideShapes.Count; $k++) {
    $cand = $slideShapes.Item($k)
    if ($cand.Name -eq "TextBox 11") {
        $ttr = $cand.TextFrame.TextRange
        $ttr.Characters(1, 16).Text = "                "
    }
}

# ---------------------------------------------------------------------
# 3) Add the "c)" and "d)" sub-figure labels (duplicates of "b)").
# ---------------------------------------------------------------------
$bBox = $null
for ($k = 1; $k -le $slideShapes.Count; $k++) {
    $cand = $slideShapes.Item($k)
    if ($cand.Name -eq "TextBox 33") {
        $bBox = $cand
    }
}

$cDup = $bBox.Duplicate()
$cBox = $cDup.Item(1)
$cBox.Name = "TextBox 16"
$cBox.TextFrame.TextRange.Characters(1, 1).Text = "c"
$cBox.Left = 3.6177165508270264
$cBox.Top = 125.81401824951172
$cBox.Width = 23.37606430053711
$cBox.Height = 21.810945510864258

$dDup = $bBox.Duplicate()
$dBox = $dDup.Item(1)
$dBox.Name = "TextBox 17"
$dBox.TextFrame.TextRange.Characters(1, 1).Text = "d"
$dBox.Left = 183.6177215576172
$dBox.Top = 125.81401824951172
$dBox.Width = 24.512126922607422
$dBox.Height = 21.810945510864258
